$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (sst index 224)
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 01:45"

$ws.Range("B4").Value = 2182717
$ws.Range("C4").Value = 20489
$ws.Range("D4").Value = 881020
$ws.Range("E4").Value = 1183418
$ws.Range("G4").Value = 421
$ws.Range("H4").Value = 118279

$ws.Range("B5").Value = 891556
$ws.Range("C5").Value = 23674
$ws.Range("E5").Value = 393870
$ws.Range("G5").Value = 729
$ws.Range("H5").Value = 44118

$ws.Range("B20").Value = 99147
$ws.Range("C20").Value = 360
$ws.Range("D20").Value = 61042
$ws.Range("E20").Value = 29930
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 8175

$ws.Range("B30").Value = 47322
$ws.Range("C30").Value = 571
$ws.Range("D30").Value = 23349
$ws.Range("E30").Value = 20044
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = 3929

$ws.Range("A37").Value = "Argentina"
$ws.Range("B37").Value = 32785
$ws.Range("C37").Value = 1208
$ws.Range("D37").Value = 9891
$ws.Range("E37").Value = 22040
$ws.Range("G37").Value = 21
$ws.Range("H37").Value = 854

$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 31810
$ws.Range("C38").Value = 656
$ws.Range("D38").Value = 14253
$ws.Range("E38").Value = 16656
$ws.Range("G38").Value = 12
$ws.Range("H38").Value = 901

$ws.Range("B52").Value = 17502
$ws.Range("C52").Value = 73
$ws.Range("D52").Value = 15686
$ws.Range("E52").Value = 891

$ws.Range("B55").Value = 16658
$ws.Range("C55").Value = 573
$ws.Range("D55").Value = 5349
$ws.Range("E55").Value = 10885
$ws.Range("G55").Value = 4
$ws.Range("H55").Value = 424

$ws.Range("B64").Value = 10064
$ws.Range("C64").Value = 40
$ws.Range("D64").Value = 7296
$ws.Range("E64").Value = 2438

$ws.Range("B69").Value = 8647
$ws.Range("C69").Value = 16
$ws.Range("E69").Value = 267

$ws.Range("A71").Value = "Sudan"
$ws.Range("B71").Value = 7435
$ws.Range("C71").Value = 215
$ws.Range("D71").Value = 2720
$ws.Range("E71").Value = 4247
$ws.Range("G71").Value = 9
$ws.Range("H71").Value = 468

$ws.Range("A72").Value = "Australia"
$ws.Range("B72").Value = 7335
$ws.Range("C72").Value = 15
$ws.Range("D72").Value = 6851
$ws.Range("E72").Value = 382
$ws.Range("H72").Value = 102

$ws.Range("B80").Value = 4572
$ws.Range("C80").Value = 40
$ws.Range("D80").Value = 3259
$ws.Range("E80").Value = 1287
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 26

$ws.Range("A86").Value = "Gabon"
$ws.Range("B86").Value = 4033
$ws.Range("C86").Value = 570
$ws.Range("D86").Value = 1334
$ws.Range("E86").Value = 2672
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 27

$ws.Range("A87").Value = "El Salvador"
$ws.Range("B87").Value = 3826
$ws.Range("C87").Value = 106
$ws.Range("D87").Value = 1912
$ws.Range("E87").Value = 1840
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 74

$ws.Range("A88").Value = "Kenia"
$ws.Range("B88").Value = 3727
$ws.Range("C88").Value = 133
$ws.Range("D88").Value = 1286
$ws.Range("E88").Value = 2337
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 104

$ws.Range("A89").Value = "Etiopia"
$ws.Range("B89").Value = 3521
$ws.Range("C89").Value = 176
$ws.Range("D89").Value = 620
$ws.Range("E89").Value = 2841
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 60

$ws.Range("A93").Value = "Venezuela"
$ws.Range("B93").Value = 3062
$ws.Range("C93").Value = 158
$ws.Range("D93").Value = 835
$ws.Range("E93").Value = 2202
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 25

$ws.Range("A94").Value = "Bosnia y Herzegovina"
$ws.Range("B94").Value = 3040
$ws.Range("C94").Value = 147
$ws.Range("D94").Value = 2162
$ws.Range("E94").Value = 713
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 165

$ws.Range("A114").Value = "Guinea-Bisau"
$ws.Range("B114").Value = 1492
$ws.Range("C114").Value = 32
$ws.Range("D114").Value = 153
$ws.Range("E114").Value = 1324
$ws.Range("H114").Value = 15

$ws.Range("A115").Value = "Libano"
$ws.Range("C115").Value = 18
$ws.Range("D115").Value = 875
$ws.Range("E115").Value = 557
$ws.Range("H115").Value = 32

$ws.Range("A116").Value = "Nicaragua"
$ws.Range("B116").Value = 1464
$ws.Range("D116").Value = 953
$ws.Range("E116").Value = 456
$ws.Range("H116").Value = 55

$ws.Range("D134").Value = 792
$ws.Range("E134").Value = 33

$ws.Range("B153").Value = 467
$ws.Range("C153").Value = 13
$ws.Range("D153").Value = 70
$ws.Range("E153").Value = 387

$ws.Range("E161").Value = 176
$ws.Range("G161").Value = 2
$ws.Range("H161").Value = 5

$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
